$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/11/2023  Through  12/17/2023"

# --- Data grid updates (rows 14-27) ---
$ws.Range("D14").NumberFormat = "#,##0"; $ws.Range("D14").Value = 1
$ws.Range("E14").NumberFormat = "#,##0.0;""-""#,##0.0"; $ws.Range("E14").Value = -100
$ws.Range("G14").NumberFormat = "#,##0"; $ws.Range("G14").Value = 1
$ws.Range("H14").NumberFormat = "#,##0.0;""-""#,##0.0"; $ws.Range("H14").Value = -100
$ws.Range("J14").NumberFormat = "#,##0"; $ws.Range("J14").Value = 1
$ws.Range("K14").NumberFormat = "#,##0.0;""-""#,##0.0"; $ws.Range("K14").Value = 0
$ws.Range("C15").NumberFormat = "@"; $ws.Range("C15").Value = "0"
$ws.Range("D15").NumberFormat = "@"; $ws.Range("D15").Value = "0"
$ws.Range("E15").NumberFormat = "@"; $ws.Range("E15").Value = "***.*"
$ws.Range("M15").Value = -12.5
$ws.Range("N15").Value = -50
$ws.Range("C16").Value = 3
$ws.Range("D16").NumberFormat = "#,##0"; $ws.Range("D16").Value = 2
$ws.Range("E16").NumberFormat = "#,##0.0;""-""#,##0.0"; $ws.Range("E16").Value = 50
$ws.Range("F16").Value = 9
$ws.Range("H16").Value = 28.571428571428
$ws.Range("I16").Value = 112
$ws.Range("J16").Value = 89
$ws.Range("K16").Value = 25.842696629213
$ws.Range("L16").Value = 23.076923076923
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = -81.487603305785
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 4
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 101
$ws.Range("J17").Value = 112
$ws.Range("K17").Value = -9.821428571428
$ws.Range("L17").Value = -14.406779661016
$ws.Range("M17").Value = 94.230769230769
$ws.Range("N17").Value = -35.668789808917
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 63.636363636363
$ws.Range("I18").Value = 126
$ws.Range("J18").Value = 178
$ws.Range("K18").Value = -29.213483146067
$ws.Range("L18").Value = -13.698630136986
$ws.Range("M18").Value = 38.461538461538
$ws.Range("N18").Value = -89.367088607594
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = -40
$ws.Range("F19").Value = 53
$ws.Range("G19").Value = 68
$ws.Range("H19").Value = -22.058823529411
$ws.Range("I19").Value = 665
$ws.Range("J19").Value = 656
$ws.Range("K19").Value = 1.371951219512
$ws.Range("L19").Value = 26.185958254269
$ws.Range("M19").Value = -7.122905027932
$ws.Range("N19").Value = -69.827586206896
$ws.Range("F20").Value = 12
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 61
$ws.Range("J20").Value = 65
$ws.Range("K20").Value = -6.153846153846
$ws.Range("L20").Value = 35.555555555555
$ws.Range("M20").Value = 103.333333333333
$ws.Range("N20").Value = -90.113452188006
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -18.75
$ws.Range("F21").Value = 97
$ws.Range("G21").Value = 101
$ws.Range("H21").Value = -3.960396039603
$ws.Range("I21").Value = 1073
$ws.Range("J21").Value = 1117
$ws.Range("K21").Value = -3.939122649955
$ws.Range("L21").Value = 14.759358288770
$ws.Range("M21").Value = 12.591815320042
$ws.Range("N21").Value = -77.561689669594
$ws.Range("C22").NumberFormat = "#,##0"; $ws.Range("C22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"; $ws.Range("D22").Value = 1
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"; $ws.Range("E22").Value = 0
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -80
$ws.Range("I22").Value = 37
$ws.Range("J22").Value = 31
$ws.Range("K22").Value = 19.354838709677
$ws.Range("L22").Value = 42.307692307692
$ws.Range("M22").Value = 8.823529411764
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 17.647058823529
$ws.Range("F24").Value = 72
$ws.Range("G24").Value = 68
$ws.Range("H24").Value = 5.882352941176
$ws.Range("I24").Value = 1024
$ws.Range("J24").Value = 1199
$ws.Range("K24").Value = -14.595496246872
$ws.Range("L24").Value = 4.918032786885
$ws.Range("M24").Value = 68.976897689769
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 16.666666666666
$ws.Range("F25").Value = 20
$ws.Range("G25").Value = 22
$ws.Range("H25").Value = -9.090909090909
$ws.Range("I25").Value = 243
$ws.Range("J25").Value = 233
$ws.Range("K25").Value = 4.291845493562
$ws.Range("L25").Value = 37.288135593220
$ws.Range("M25").Value = 9.459459459459
$ws.Range("C26").NumberFormat = "@"; $ws.Range("C26").Value = "0"
$ws.Range("D26").NumberFormat = "@"; $ws.Range("D26").Value = "0"
$ws.Range("E26").NumberFormat = "@"; $ws.Range("E26").Value = "***.*"
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -75
$ws.Range("J27").Value = 54
$ws.Range("K27").Value = 16.666666666666
$ws.Range("L27").Value = -8.695652173913
